$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'AU-4,AU-4 (1)'
$ws.Cells.Item(3, 1).Value = 'SC-5,CM-6 b,SC-5 (2)'
$ws.Cells.Item(4, 1).Value = 'AU-7 a,AC-6 (8),AU-7 b,AU-12 (3),CM-5 (1),AC-6 (9),AU-8 b'
$ws.Cells.Item(5, 1).Value = 'CM-7 b,AC-17 (9),AC-17 (1),CM-6 b'
$ws.Cells.Item(10, 1).Value = 'CM-7 (2),CM-7 (5) (b)'
$ws.Cells.Item(17, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(19, 1).Value = 'IA-5 (1) (a),IA-5 (1) (b),CM-6 b'
$ws.Cells.Item(22, 1).Value = 'AU-7 a,AU-3 (1),CM-6 b,AU-6 (4),CM-5 (1),AU-7 (1),MA-4 (1) (a),AU-14 (1),AU-3,AU-12 a'
$ws.Cells.Item(25, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(29, 1).Value = 'SC-8,SC-8 (1),SC-8 (2)'
$ws.Cells.Item(31, 1).Value = 'AU-12 c,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(38, 1).Value = 'AU-9,SI-11 b'
$ws.Cells.Item(42, 1).Value = 'SC-28,SC-28 (1)'
$ws.Cells.Item(45, 1).Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a,AC-8 b'
$ws.Cells.Item(50, 1).Value = 'IA-2 (5),CM-6 b'
$ws.Cells.Item(53, 1).Value = 'SC-13,MA-4 (6)'
$ws.Cells.Item(55, 1).Value = 'SC-8,AC-17 (2)'
$ws.Cells.Item(56, 1).Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Cells.Item(63, 1).Value = 'AU-5 (1),AU-5 a'
$ws.Cells.Item(67, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(69, 1).Value = 'AU-12 c,AU-7 a,CM-6 b,AU-12 (3),AU-7 b,CM-5 (1),AU-8 b,AU-12 a'
$ws.Cells.Item(77, 1).Value = 'AU-12 c,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(80, 1).Value = 'IA-2 (2),IA-2 (1),IA-2 (3),IA-2 (4)'
$ws.Cells.Item(86, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(88, 1).Value = 'AU-12 c,CM-5 (1),AC-6 (9),AC-2 (4)'
$ws.Cells.Item(89, 1).Value = 'IA-2 (5),IA-2 (3),IA-2 (4),IA-2,IA-2 (2)'
$ws.Cells.Item(96, 1).Value = 'SC-8,SC-8 (1),AC-18 (1)'
$ws.Cells.Item(97, 1).Value = 'AU-8 (1) (b),AU-8 b,AU-8 (1) (a)'
$ws.Cells.Item(102, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(111, 1).Value = 'AU-5 b,AU-5 a'
$ws.Cells.Item(119, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(123, 1).Value = 'CM-7 b,CM-7 a'
$ws.Cells.Item(124, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(125, 1).Value = 'AC-18 (1),CM-7 a'
$ws.Cells.Item(128, 1).Value = 'CM-7 a,CM-6 b,IA-5 (1) (c)'
$ws.Cells.Item(148, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-14 (1),AU-3,AU-12 a'
$ws.Cells.Item(157, 1).Value = 'AU-12 c,AU-3 (1),MA-4 (1) (a),AU-3,AU-12 a'
$ws.Cells.Item(159, 1).Value = 'SC-8,AC-17 (2)'
$ws.Cells.Item(175, 1).Value = 'SI-16,CM-7 a'
